$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the Price (D) and Volume (E) columns for the rows
# being updated, so values like "56.928.28" or "0.999" and percentage strings
# such as "  +4.33%  " are written as literal text (matching the source data,
# which stores these columns as inline strings), then restore the default
# "Normal" style so no stray style index gets attached to the cells.

$rng = $ws.Range("D2","E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = '56.928.28'
$ws.Range("E2").Value = '  +4.33%  '
$rng.Style = "Normal"

$rng = $ws.Range("D3","E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = '3.243.37'
$ws.Range("E3").Value = '  +1.80%  '
$rng.Style = "Normal"

$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '
$rng.Style = "Normal"

$rng = $ws.Range("D5","E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = '395.16'
$ws.Range("E5").Value = '  -1.95%  '
$rng.Style = "Normal"

$rng = $ws.Range("D6","E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = '108.14'
$ws.Range("E6").Value = '  -0.95%  '
$rng.Style = "Normal"

$rng = $ws.Range("D7","E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = '0.583'
$ws.Range("E7").Value = '  +5.79%  '
$rng.Style = "Normal"

$rng = $ws.Range("D8","E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = '3.239.41'
$ws.Range("E8").Value = '  +1.89%  '
$rng.Style = "Normal"

$rng = $ws.Range("D9","E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.03%  '
$rng.Style = "Normal"

$rng = $ws.Range("D10","E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = '0.622'
$ws.Range("E10").Value = '  +0.16%  '
$rng.Style = "Normal"

$rng = $ws.Range("D11","E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = '39.05'
$ws.Range("E11").Value = '  +0.26%  '
$rng.Style = "Normal"

$rng = $ws.Range("D12","E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = '0.0992'
$ws.Range("E12").Value = '  +12.31%  '
$rng.Style = "Normal"

$rng = $ws.Range("E13")
$rng.NumberFormat = "@"
$ws.Range("E13").Value = '  +1.70%  '
$rng.Style = "Normal"

$rng = $ws.Range("D14","E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = '3.758.30'
$ws.Range("E14").Value = '  +2.24%  '
$rng.Style = "Normal"

$rng = $ws.Range("D15","E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = '8.27'
$ws.Range("E15").Value = '  +2.60%  '
$rng.Style = "Normal"

$rng = $ws.Range("D16","E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = '19.06'
$ws.Range("E16").Value = '  -0.90%  '
$rng.Style = "Normal"

$rng = $ws.Range("D17","E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = '3.263.86'
$ws.Range("E17").Value = '  +2.87%  '
$rng.Style = "Normal"

$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$ws.Range("E18").Value = '  -3.40%  '
$rng.Style = "Normal"

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$ws.Range("E19").Value = '  +2.02%  '
$rng.Style = "Normal"

$rng = $ws.Range("D20","E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = '56.758.80'
$ws.Range("E20").Value = '  +4.31%  '
$rng.Style = "Normal"

$rng = $ws.Range("D21","E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = '3.34'
$ws.Range("E21").Value = '  +1.12%  '
$rng.Style = "Normal"

$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$ws.Range("E22").Value = '  +10.23%  '
$rng.Style = "Normal"

$rng = $ws.Range("D23","E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = '12.87'
$ws.Range("E23").Value = '  +0.10%  '
$rng.Style = "Normal"

$rng = $ws.Range("D24","E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = '295.35'
$ws.Range("E24").Value = '  +7.64%  '
$rng.Style = "Normal"

$rng = $ws.Range("D25","E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = '74.19'
$ws.Range("E25").Value = '  +3.47%  '
$rng.Style = "Normal"

$rng = $ws.Range("E26")
$rng.NumberFormat = "@"
$ws.Range("E26").Value = '  -3.68%  '
$rng.Style = "Normal"

$rng = $ws.Range("D27","E27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = '28.06'
$ws.Range("E27").Value = '  +1.09%  '
$rng.Style = "Normal"

$rng = $ws.Range("D28","E28")
$rng.NumberFormat = "@"
$ws.Range("D28").Value = '4.36'
$ws.Range("E28").Value = '  +0.99%  '
$rng.Style = "Normal"

$rng = $ws.Range("D29","E29")
$rng.NumberFormat = "@"
$ws.Range("D29").Value = '7.59'
$ws.Range("E29").Value = '  -5.07%  '
$rng.Style = "Normal"

$rng = $ws.Range("D30","E30")
$rng.NumberFormat = "@"
$ws.Range("D30").Value = '7.23'
$ws.Range("E30").Value = '  -2.20%  '
$rng.Style = "Normal"

$rng = $ws.Range("E31")
$rng.NumberFormat = "@"
$ws.Range("E31").Value = '  -1.57%  '
$rng.Style = "Normal"

$rng = $ws.Range("E32")
$rng.NumberFormat = "@"
$ws.Range("E32").Value = '  -0.03%  '
$rng.Style = "Normal"

$rng = $ws.Range("D33","E33")
$rng.NumberFormat = "@"
$ws.Range("D33").Value = '11.26'
$ws.Range("E33").Value = '  +1.59%  '
$rng.Style = "Normal"

$rng = $ws.Range("E34")
$rng.NumberFormat = "@"
$ws.Range("E34").Value = '  -4.02%  '
$rng.Style = "Normal"

$rng = $ws.Range("D35","E35")
$rng.NumberFormat = "@"
$ws.Range("D35").Value = '39.21'
$ws.Range("E35").Value = '  +6.69%  '
$rng.Style = "Normal"

$rng = $ws.Range("D36","E36")
$rng.NumberFormat = "@"
$ws.Range("D36").Value = '0.0480'
$ws.Range("E36").Value = '  -3.81%  '
$rng.Style = "Normal"

$rng = $ws.Range("E37")
$rng.NumberFormat = "@"
$ws.Range("E37").Value = '  +2.11%  '
$rng.Style = "Normal"

$rng = $ws.Range("D38","E38")
$rng.NumberFormat = "@"
$ws.Range("D38").Value = '51.42'
$ws.Range("E38").Value = '  +1.31%  '
$rng.Style = "Normal"

$rng = $ws.Range("D39","E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.02%  '
$rng.Style = "Normal"

$rng = $ws.Range("E40")
$rng.NumberFormat = "@"
$ws.Range("E40").Value = '  -5.70%  '
$rng.Style = "Normal"

$rng = $ws.Range("E41")
$rng.NumberFormat = "@"
$ws.Range("E41").Value = '  +1.65%  '
$rng.Style = "Normal"

$rng = $ws.Range("D42","E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = '135.07'
$ws.Range("E42").Value = '  +3.73%  '
$rng.Style = "Normal"

$rng = $ws.Range("E43")
$rng.NumberFormat = "@"
$ws.Range("E43").Value = '  +4.03%  '
$rng.Style = "Normal"

$rng = $ws.Range("E44")
$rng.NumberFormat = "@"
$ws.Range("E44").Value = '  -2.59%  '
$rng.Style = "Normal"

$rng = $ws.Range("D45","E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = '3.94'
$ws.Range("E45").Value = '  -3.90%  '
$rng.Style = "Normal"

$rng = $ws.Range("D46","E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = '16.98'
$ws.Range("E46").Value = '  -1.93%  '
$rng.Style = "Normal"

$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$ws.Range("E47").Value = '  -4.42%  '
$rng.Style = "Normal"

$rng = $ws.Range("D48","E48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = '22.06'
$ws.Range("E48").Value = '  -0.99%  '
$rng.Style = "Normal"

$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$ws.Range("E49").Value = '  +3.13%  '
$rng.Style = "Normal"

$rng = $ws.Range("D50","E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = '2.148.27'
$ws.Range("E50").Value = '  +2.87%  '
$rng.Style = "Normal"

$rng = $ws.Range("D51","E51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = '2.34'
$ws.Range("E51").Value = '  -6.91%  '
$rng.Style = "Normal"
